# NIT-9012257682.xlsx - "Actualiza base de datos EC y agrega parte 1 de
# nuevos estado de cuenta"
#
# This updates the EC (estado de cuenta) worksheet:
#  - The "Valor Mora" total (E11) and the worker/period counts (C13, F13)
#    are refreshed.
#  - The old detail row for worker MARIANNY JOSE MOSQUERA PADILLA with the
#    stale document number (5295483 / period 2312) is removed; the
#    up-to-date record for the same worker (document 20475746 / period
#    2401) takes its place, shifting the signature block up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the summary figures.
$ws.Range("E11").Value = 14142   # VALOR MORA
$ws.Range("C13").Value = 2       # Cant. Trabajadores
$ws.Range("F13").Value = 2       # Cant. Periodos

# Drop the outdated MARIANNY JOSE MOSQUERA PADILLA row (doc 5295483,
# period 2312). The row below it (doc 20475746, period 2401) slides up
# into its place, and the trailing signature rows move up with it.
$ws.Rows(17).Delete()
